$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 2168.19025045373
$ws.Range("C2").Value = 1066.827089534028
$ws.Range("D2").Value = 1207.57531266709
$ws.Range("E2").Value = 2196.599662171774
$ws.Range("F2").Value = 2169.953739032937
$ws.Range("G2").Value = 1810.682651439943
$ws.Range("H2").Value = 2197.755167519784

$ws.Range("B3").Value = 2183.438124058833
$ws.Range("C3").Value = 1096.168533900459
$ws.Range("D3").Value = 1063.261003313442
$ws.Range("E3").Value = 2201.606885761249
$ws.Range("F3").Value = 2184.514313820853
$ws.Range("G3").Value = 1768.334091166898
$ws.Range("H3").Value = 2202.387130431995

$ws.Range("B4").Value = 2134.223707206921
$ws.Range("C4").Value = 1121.183306044844
$ws.Range("D4").Value = 1046.252314260004
$ws.Range("E4").Value = 2178.408688171938
$ws.Range("F4").Value = 2137.004174443122
$ws.Range("G4").Value = 1771.385937107427
$ws.Range("H4").Value = 2180.219103371565

$ws.Range("B5").Value = 2184.588446198372
$ws.Range("C5").Value = 1125.261196947115
$ws.Range("D5").Value = 1050.716306735508
$ws.Range("E5").Value = 2198.951043606867
$ws.Range("F5").Value = 2185.240087232406
$ws.Range("G5").Value = 1696.724731045043
$ws.Range("H5").Value = 2199.509380308793

$ws.Range("B6").Value = 2190.810902805146
$ws.Range("C6").Value = 1103.837453310326
$ws.Range("D6").Value = 362.6284403737303
$ws.Range("E6").Value = 2206.670585597019
$ws.Range("F6").Value = 2191.438427965602
$ws.Range("G6").Value = 1252.015568647364
$ws.Range("H6").Value = 2207.103148948721

$ws.Range("B7").Value = 2172.524713790897
$ws.Range("C7").Value = 1154.182764083716
$ws.Range("D7").Value = 535.7477298025481
$ws.Range("E7").Value = 2191.721997665262
$ws.Range("F7").Value = 2173.658829612546
$ws.Range("G7").Value = 1349.580071416546
$ws.Range("H7").Value = 2192.770205129404

$ws.Range("B8").Value = 2174.232017974375
$ws.Range("C8").Value = 1048.501688302613
$ws.Range("D8").Value = 633.0908430649812
$ws.Range("E8").Value = 2188.904419689039
$ws.Range("F8").Value = 2174.829099415711
$ws.Range("G8").Value = 1347.754846155285
$ws.Range("H8").Value = 2189.482557421556

$ws.Range("B9").Value = 2186.103315481719
$ws.Range("C9").Value = 1134.340219809489
$ws.Range("D9").Value = 616.5541332604442
$ws.Range("E9").Value = 2194.535756323377
$ws.Range("F9").Value = 2187.41046822575
$ws.Range("G9").Value = 1364.57542459144
$ws.Range("H9").Value = 2195.799347806329

$ws.Range("B10").Value = 1946.701420854897
$ws.Range("C10").Value = 1228.166133370289
$ws.Range("D10").Value = 1004.768125513142
$ws.Range("E10").Value = 2045.207259279136
$ws.Range("F10").Value = 1963.41524271665
$ws.Range("G10").Value = 1657.608033392178
$ws.Range("H10").Value = 2056.801072003444

$ws.Range("B11").Value = 1895.566769401754
$ws.Range("C11").Value = 1294.365777851517
$ws.Range("D11").Value = 843.8725314057115
$ws.Range("E11").Value = 2033.936097627198
$ws.Range("F11").Value = 1912.761164030131
$ws.Range("G11").Value = 1635.170310488515
$ws.Range("H11").Value = 2044.274956698688

$ws.Range("B12").Value = 1637.698552672618
$ws.Range("C12").Value = 1172.687020782348
$ws.Range("D12").Value = 198.0395115404653
$ws.Range("E12").Value = 1876.610898948367
$ws.Range("F12").Value = 1646.942019134947
$ws.Range("G12").Value = 1247.307099901423
$ws.Range("H12").Value = 1880.46872151288

$ws.Range("B13").Value = 1929.418281600073
$ws.Range("C13").Value = 1219.485333534142
$ws.Range("D13").Value = 752.61671157398
$ws.Range("E13").Value = 2034.261254671876
$ws.Range("F13").Value = 1942.524931502256
$ws.Range("G13").Value = 1529.951663943509
$ws.Range("H13").Value = 2043.082140312506
